$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab: "Scanner" -> "Anatomy_Lecture_Hall"
$ws.Name = "Anatomy_Lecture_Hall"

# Header row: E1 "Type" -> "Number"
$ws.Range("E1").Value = "Number"

# Data row 2 updates.
# C2 and E2 look like a date and a big integer respectively, so Excel would
# normally auto-convert them on assignment. Force those cells to Text format
# first so the values are stored verbatim as strings, matching the source data.
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "04/05/2025"

# D2 "3:45 PM" -> "15:45:09" (not auto-parsed as a time by Excel, stays text)
$ws.Range("D2").Value = "15:45:09"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1746362709541"
